$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.362.37"
$ws.Range("E2").Value = "  +1.26%  "

$ws.Range("D3").Value = "'1.685.21"
$ws.Range("E3").Value = "  +1.03%  "

$ws.Range("D4").Value = "'1.007"
$ws.Range("E4").Value = "  +0.24%  "

$ws.Range("D5").Value = "'218.48"
$ws.Range("E5").Value = "  +0.73%  "

$ws.Range("D6").Value = "'0.5578"
$ws.Range("E6").Value = "  +9.24%  "

$ws.Range("D7").Value = "'1.007"
$ws.Range("E7").Value = "  +0.25%  "

$ws.Range("D8").Value = "'0.2707"
$ws.Range("E8").Value = "  +1.78%  "

$ws.Range("D9").Value = "'0.06505"
$ws.Range("E9").Value = "  +1.59%  "

$ws.Range("D10").Value = "'22.17"
$ws.Range("E10").Value = "  +1.78%  "

$ws.Range("D11").Value = "'0.07564"
$ws.Range("E11").Value = "  +1.62%  "

$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").Value = "'4.550"
$ws.Range("E12").Value = "  +0.78%  "

$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "'1.685.25"
$ws.Range("E13").Value = "  +0.89%  "

$ws.Range("D14").Value = "'0.5812"
$ws.Range("E14").Value = "  -0.37%  "

$ws.Range("D15").Value = "'0.000008462"
$ws.Range("E15").Value = "  -1.05%  "

$ws.Range("D16").Value = "'65.40"
$ws.Range("E16").Value = "  +1.72%  "

$ws.Range("D17").Value = "'26.395.51"
$ws.Range("E17").Value = "  +1.17%  "

$ws.Range("D18").Value = "'4.947"
$ws.Range("E18").Value = "  +0.11%  "

$ws.Range("D19").Value = "'1.007"
$ws.Range("E19").Value = "  +0.27%  "

$ws.Range("D20").Value = "'10.93"
$ws.Range("E20").Value = "  +1.40%  "

$ws.Range("D21").Value = "'191.53"
$ws.Range("E21").Value = "  -0.15%  "

$ws.Range("D22").Value = "'6.246"
$ws.Range("E22").Value = "  +0.73%  "

$ws.Range("D23").Value = "'1.007"
$ws.Range("E23").Value = "  +0.16%  "

$ws.Range("D24").Value = "'148.09"
$ws.Range("E24").Value = "  +2.27%  "

$ws.Range("D25").Value = "'0.1333"
$ws.Range("E25").Value = "  +11.04%  "

$ws.Range("D26").Value = "'7.910"
$ws.Range("E26").Value = "  +3.87%  "

$ws.Range("D27").Value = "'15.84"
$ws.Range("E27").Value = "  +1.10%  "

$ws.Range("D28").Value = "'0.06334"
$ws.Range("E28").Value = "  -2.62%  "

$ws.Range("D29").Value = "'1.394"
$ws.Range("E29").Value = "  +4.13%  "

$ws.Range("E30").Value = "  +0.56%  "

$ws.Range("D31").Value = "'3.598"
$ws.Range("E31").Value = "  +1.62%  "

$ws.Range("D32").Value = "'3.584"
$ws.Range("E32").Value = "  +1.94%  "

$ws.Range("D33").Value = "'1.673"
$ws.Range("E33").Value = "  +1.33%  "

$ws.Range("D34").Value = "'1.044"
$ws.Range("E34").Value = "  +2.39%  "

$ws.Range("D35").Value = "'0.6218"
$ws.Range("E35").Value = "  +1.66%  "

$ws.Range("D36").Value = "'2.400"
$ws.Range("E36").Value = "  +1.23%  "

$ws.Range("D37").Value = "'2.721"
$ws.Range("E37").Value = "  +1.38%  "

$ws.Range("D38").Value = "'6.237"
$ws.Range("E38").Value = "  -0.43%  "

$ws.Range("D39").Value = "'1.114.79"
$ws.Range("E39").Value = "  +2.07%  "

$ws.Range("D40").Value = "'0.01627"
$ws.Range("E40").Value = "  +1.59%  "

$ws.Range("D41").Value = "'0.8766"
$ws.Range("E41").Value = "  +0.92%  "

$ws.Range("D42").Value = "'1.014"
$ws.Range("E42").Value = "  +0.46%  "

$ws.Range("D43").Value = "'100.62"
$ws.Range("E43").Value = "  -0.63%  "

$ws.Range("D44").Value = "'1.834.22"

$ws.Range("E45").Value = "  -4.53%  "

$ws.Range("D46").Value = "'57.44"
$ws.Range("E46").Value = "  +1.74%  "

$ws.Range("D47").Value = "'8.218"
$ws.Range("E47").Value = "  +1.65%  "

$ws.Range("D48").Value = "'1.005"
$ws.Range("E48").Value = "  -0.01%  "

$ws.Range("E49").Value = "  +0.81%  "

$ws.Range("D50").Value = "'0.4296"
$ws.Range("E50").Value = "  +0.23%  "

$ws.Range("D51").Value = "'6.089"
$ws.Range("E51").Value = "  -0.02%  "
